$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "SEQ_NO(M)"
$ws.Range("B1").Value = "HOUSE_BL_NO(M)"
$ws.Range("C1").Value = "SHIPMENT_DATE(YYYYMMDD)(M)"
$ws.Range("D1").Value = "ARRIVAL_DATE(YYYYMMDD)(M)"
$ws.Range("E1").Value = "SHIPMENT_COUNTRY_CODE(M)"
$ws.Range("F1").Value = "SHIPMENT_METHOD_CODE(M)"
$ws.Range("G1").Value = "ARRIVAL_PORT_CODE(M)"
$ws.Range("H1").Value = "PRICE_TERMS_CODE(M)"
$ws.Range("I1").Value = "SHIPMENT_CURRENCY_CODE(M)"
$ws.Range("J1").Value = "FINAL_DESTINATION_CODE(M)"
$ws.Range("K1").Value = "LINE_NUM(M)"
$ws.Range("L1").Value = "ITEM_NO(M)"
$ws.Range("M1").Value = "QUANTITY_SHIPPED(M)"
$ws.Range("N1").Value = "UNIT_PRICE(M)"
$ws.Range("O1").Value = "VAT_AMOUNT"
$ws.Range("P1").Value = "IMPORT_DECLARATION_NO"
$ws.Range("Q1").Value = "IMPORT_DECLARATION_DATE(YYYYMMDD)"
$ws.Range("R1").Value = "CC_AMOUNT"
$ws.Range("S1").Value = "DUTY_AMOUNT"
$ws.Range("T1").Value = "SHIPMENT_PORT_CODE"
$ws.Range("U1").Value = "PACKING_QTY"
$ws.Range("V1").Value = "PACKING_UNIT_CODE"
$ws.Range("W1").Value = "CARGO_TYPE_CODE"
$ws.Range("X1").Value = "TOTAL_CBM"
$ws.Range("Y1").Value = "GROSS_WEIGHT"
$ws.Range("Z1").Value = "CHARGEABLE_WEIGHT"
$ws.Range("AA1").Value = "GROSS_WEIGHT_UOM_CODE"
$ws.Range("AB1").Value = "VESSEL_FLIGHT_NAME"
$ws.Range("AC1").Value = "VENDOR_SITE_CODE"
$ws.Range("AD1").Value = "COURIER_FLAG"
$ws.Range("AE1").Value = "ITEM_SPEC"
$ws.Range("AF1").Value = "HS_NO"
$ws.Range("AG1").Value = "UNIT_OF_MEASURE"
$ws.Range("AH1").Value = "CC_COUNTRY_CODE"
$ws.Range("AI1").Value = "ORGANIZATION_CODE"

$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "2395544244"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "KR"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "DEMUC"
$ws.Range("H2").Value = "DAP"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = "DEMUC"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = "LGE CID G01 10.25 APX2T FST"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "996,68"
$ws.Range("P2").Value = "ATC00400931750820195604"
$ws.Range("Q2").Value = "20190819"
$ws.Range("R2").Value = "5.117,76"
$ws.Range("S2").Value = "127,94"
$ws.Range("T2").Value = ""
$ws.Range("U2").Value = ""
$ws.Range("V2").Value = ""
$ws.Range("W2").Value = ""
$ws.Range("X2").Value = ""
$ws.Range("Y2").Value = ""
$ws.Range("Z2").Value = ""
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = ""
$ws.Range("AC2").Value = ""
$ws.Range("AD2").Value = ""
$ws.Range("AE2").Value = ""
$ws.Range("AF2").Value = "85299092990"
$ws.Range("AG2").Value = ""
$ws.Range("AH2").Value = ""
$ws.Range("AI2").Value = ""